# Auto-committed on 2022/08/26 週五 17:20:49.47
#
# Adds two new field definitions (rows 16 & 17 of the DBD table layout)
# for JcicZ048: "ActualFilingDate" (實際報送日期) and
# "ActualFilingMark" (實際報送記號), and updates the saved selection on
# the DBD sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 24 -> field #16 : ActualFilingDate
$ws.Range("A24").Value = 16
$ws.Range("B24").Value = "ActualFilingDate"
$ws.Range("C24").Value = "實際報送日期"
$ws.Range("D24").Value = "Decimald"
$ws.Range("E24").Value = 8

# Row 25 -> field #17 : ActualFilingMark
$ws.Range("A25").Value = 17
$ws.Range("B25").Value = "ActualFilingMark"
$ws.Range("C25").Value = "實際報送記號"
$ws.Range("D25").Value = "VARCHAR2"
$ws.Range("E25").Value = 3

# Restore the sheet as active and move the saved selection/scroll
# position to where the author left it (B28) after typing the new rows.
$ws.Activate() | Out-Null
$ws.Range("B28").Select() | Out-Null
